# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on row 14 of the
# zh-cn and de-de sheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("D14").Value = "2016-03-02 15:29:07"
$wsZhCn.Range("G14").Value = "2016-03-02 15:29:50"

$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("D14").Value = "2016-03-02 15:29:18"
$wsDeDe.Range("G14").Value = "2016-03-02 15:30:12"
